$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark (Word always keeps this bookmark at the
#    location of the most recent edit, so it needs to move to the new text).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Find the empty paragraph right after the "Source Control" heading and
#    fill it in with the new write-up, then restore a trailing empty
#    paragraph (mirroring the one that used to be there) after it.
# ---------------------------------------------------------------------------
$headingIdx = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text.Trim() -eq "Source Control") {
        $headingIdx = $idx
    }
}
$target = $d.Paragraphs.Item($headingIdx + 1)

# Normalise the (still empty) paragraph mark formatting: it used to carry the
# bold / 36pt "heading" look (inherited so a *new* heading could be typed),
# but now that it holds body text it should fall back to plain formatting.
$target.Range.InsertAfter("X")
$markRange = $target.Range
$markRange.Collapse(0)
$markRange.Font.Bold = $false
$markRange.Font.Size = 11
$placeholder = $d.Range($target.Range.Start, $target.Range.Start + 1)
$placeholder.Delete()

# Give the paragraph some space after it, as in the final layout.
$target.Format.SpaceAfter = 12

# Add a new, empty paragraph after this one - it keeps the bold/36pt look
# that used to live on the paragraph we just repurposed.
$target.Range.InsertParagraphAfter()

# Type in all of the new text first (so the insertion point's "current
# format" never gets contaminated by a mid-stream Font change), remembering
# each chunk's [start,end) so the sz=24 runs can be formatted afterwards.
$c0s = $target.Range.End - 1
$target.Range.InsertAfter("In this project we had a chance to use the GitHub as we have to. We create a repository which called myTunes and also we attached its link to the frontpage.  In this repo")
$c0e = $target.Range.End - 1

$c1s = $target.Range.End - 1
$target.Range.InsertAfter("sitory we can find the source code, the image and a .zip file which contains all the java libraries what we used. During the project it turned out that the GitHub is a well-constructed software and it contributed to our colla")
$c1e = $target.Range.End - 1

$target.Range.InsertAfter("borative work")
$target.Range.InsertAfter(". ")
$target.Range.InsertAfter("Even though")
$target.Range.InsertAfter(" we experienced the benefits of GitHub, we ")
$target.Range.InsertAfter("had an opportunity to get to know its drawback. This drawback is the merge conflict what we didn’t know before. As we know it occurs when the same part of the code exit in two different way but the same name. We had to solve all the merge conflicts during the project which improved our problem-solving skills. ")

# Now go back and bump the first two chunks up to 12pt (sz=24).
$chunk0 = $d.Range($c0s, $c0e)
$chunk0.Font.Size = 12
$chunk1 = $d.Range($c1s, $c1e)
$chunk1.Font.Size = 12

# Re-insert the "_GoBack" bookmark at the very end of the new text (Word
# drops it at the last edited location).
$endOfParagraph = $d.Range($target.Range.End - 1, $target.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $endOfParagraph)

Write-Output "done"
